$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the numeric-iteration result rows (2-4) with the new computed values.
# NumberFormat is forced to text ("@") first so the numeric-looking strings
# (e.g. "-2.0", "nan") are stored as text, matching the original inlineStr data.

$ws.Range("B2:D4").NumberFormat = "@"

$ws.Range("B2").Value = "-2.0"
$ws.Range("C2").Value = "-8.0"
$ws.Range("D2").Value = "2.0003606680642"

$ws.Range("B3").Value = "0.000360668064204"
$ws.Range("C3").Value = "4.6916225675279003e-11"
$ws.Range("D3").Value = "0.0003576882720103"

$ws.Range("B4").Value = "2.97979219374231e-06"
$ws.Range("C4").Value = "2.64580561779787e-17"
$ws.Range("D4").Value = "nan"

# Remove the now-obsolete trailing iteration rows (5-8), shrinking the
# worksheet's used range/dimension from A1:D8 down to A1:D4.
$ws.Range("A5:D8").EntireRow.Delete()
